$wb = $excel.ActiveWorkbook

# --- Sheet "Means" ---
$ws = $wb.Worksheets.Item("Means")

# Row 9 (Total Cancer Risk (per million))
$ws.Range("B9").Value = 26
$ws.Range("C9").Value = 26
$ws.Range("F9").Value = 31

# Row 10 (Total Respiratory (hazard quotient))
$ws.Range("B10").Value = 0.31
$ws.Range("C10").Value = 0.32
$ws.Range("D10").Value = 0.36
$ws.Range("E10").Value = 0.37
$ws.Range("F10").Value = 0.38
$ws.Range("G10").Value = 0.36

# --- Sheet "Standard Deviations" ---
$ws2 = $wb.Worksheets.Item("Standard Deviations")

# Row 9 (Total Cancer Risk (per million) SD)
$ws2.Range("B9").Value = 8.3
$ws2.Range("C9").Value = 5.4
$ws2.Range("F9").Value = 6.4
$ws2.Range("G9").Value = 4.2

# Row 10 (Total Respiratory (hazard quotient) SD)
$ws2.Range("B10").Value = 0.11
$ws2.Range("C10").Value = 0.043
$ws2.Range("D10").Value = 0.055
$ws2.Range("E10").Value = 0.044
$ws2.Range("F10").Value = 0.048
$ws2.Range("G10").Value = 0.052
